$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "situacao" column (J) values: each "meses sem comprar" counter
# advanced by 0.1 of a month for these clients.
$ws.Range("J9").Value  = "INATIVO - 18.1 meses sem comprar"
$ws.Range("J16").Value = "INATIVO - 1.5 meses sem comprar"
$ws.Range("J22").Value = "INATIVO - 28.2 meses sem comprar"
$ws.Range("J31").Value = "INATIVO - 6.6 meses sem comprar"
$ws.Range("J34").Value = "INATIVO - 25.8 meses sem comprar"
$ws.Range("J41").Value = "INATIVO - 6.5 meses sem comprar"
$ws.Range("J42").Value = "INATIVO - 24.8 meses sem comprar"
$ws.Range("J49").Value = "INATIVO - 9.0 meses sem comprar"
$ws.Range("J66").Value = "INATIVO - 11.6 meses sem comprar"
$ws.Range("J67").Value = "INATIVO - 10.4 meses sem comprar"
$ws.Range("J79").Value = "INATIVO - 21.3 meses sem comprar"
$ws.Range("J83").Value = "INATIVO - 20.8 meses sem comprar"
$ws.Range("J84").Value = "INATIVO - 8.5 meses sem comprar"
$ws.Range("J98").Value = "INATIVO - 21.9 meses sem comprar"
$ws.Range("J107").Value = "INATIVO - 20.4 meses sem comprar"

# Update row 111 (BEMOL S/A): new purchase-history total and refreshed
# ultima_compra / proxima_compra timestamps (roughly +1 day, new time-of-day).
$ws.Range("E111").Value = 14998
$ws.Range("H111").Value = (Get-Date -Year 2025 -Month 5 -Day 29 -Hour 17 -Minute 21 -Second 34)
$ws.Range("I111").Value = (Get-Date -Year 2025 -Month 5 -Day 30 -Hour 17 -Minute 21 -Second 34)
